$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("web 120")

# Update O13 to use a formula instead of a literal 0
$ws.Range("O13").Formula = "=(40/40)*10"

# Update Project (D column) values for several students
$ws.Range("D24").Value = 55
$ws.Range("D25").Value = 50
$ws.Range("D26").Value = 60
$ws.Range("D27").Value = 52
$ws.Range("D31").Value = 32
$ws.Range("D32").Value = 55
$ws.Range("D33").Value = 50
$ws.Range("D36").Value = 55

# Update the selected cell to match the saved view state
$ws.Range("H32").Select()

$wb.Save()
